$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 previously held the "Coordinator" shift; it becomes "Logistics".
$ws.Range("D5").Value = "Logistics"
$ws.Range("A5").Value = "TurnoLog"

# New row 6: Business shift.
$ws.Range("A6").Value = "TurnoBus"
$ws.Range("D6").Value = "Business"
$ws.Range("C6").Value = 2

# The "decorre" availability pattern grew from 10 to 20 entries; update it
# for every data row (2-6).
$pattern20 = "[1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1]"
$ws.Range("B2").Value = $pattern20
$ws.Range("B3").Value = $pattern20
$ws.Range("B4").Value = $pattern20
$ws.Range("B5").Value = $pattern20
$ws.Range("B6").Value = $pattern20

# Update the selection to match the active cell of the edit.
$ws.Range("B6").Select()
